$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This change regenerates the handback-status report for a new run:
#   old handoff file "25699400-7c3f-4264-a96c-bd7921a169d5.md"
#     -> new handoff file "836d1c83-52d1-4579-9d23-7f3bdff4659d.md"
#   old handoff file "da1b3d5f-85d4-4d96-bfee-43fbf7b22681.md"
#     -> new handoff file "ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md"
# along with refreshed xliff names and timestamps.
# ---------------------------------------------------------------------------

# ===================== Sheet 1: Overview =====================
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.md"
$ws1.Range("B2").Value = "e2e\836d1c83-52d1-4579-9d23-7f3bdff4659d.md"
$ws1.Range("G2").Value = "2016-09-01 09:20:15"

$ws1.Range("A3").Value = "ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md"
$ws1.Range("B3").Value = "e2e\ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md"
$ws1.Range("G3").Value = "2016-09-01 09:20:15"

# Hyperlinks keep the same underlying target (rId2 / rId3 URLs stay the
# same), only the displayed text changes. The engine's Hyperlink.TextToDisplay
# setter creates a duplicate entry instead of updating in place, so rebuild
# the collection: clear it, then re-add with the original target URLs.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f696022f6e3b5d387715a091f95e4dadf4e46478/e2e/25699400-7c3f-4264-a96c-bd7921a169d5.md", "", "", "e2e\836d1c83-52d1-4579-9d23-7f3bdff4659d.md")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f696022f6e3b5d387715a091f95e4dadf4e46478/e2e/da1b3d5f-85d4-4d96-bfee-43fbf7b22681.md", "", "", "e2e\ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md")

# ===================== Sheet 2: zh-cn =====================
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.md"
$ws2.Range("G2").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.03b0573f4c3ad94c7c118a0f76853d51b4646da2.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-09-01 09:19:58"
$ws2.Range("I2").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.md"
$ws2.Range("J2").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.03b0573f4c3ad94c7c118a0f76853d51b4646da2.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-09-01 09:20:36"

$ws2.Range("A3").Value = "ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md"
$ws2.Range("G3").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.03b0573f4c3ad94c7c118a0f76853d51b4646da2.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-09-01 09:19:58"
$ws2.Range("I3").Value = "ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md"
$ws2.Range("J3").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.03b0573f4c3ad94c7c118a0f76853d51b4646da2.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-09-01 09:20:36"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f696022f6e3b5d387715a091f95e4dadf4e46478/e2e/25699400-7c3f-4264-a96c-bd7921a169d5.md", "", "", "836d1c83-52d1-4579-9d23-7f3bdff4659d.md")
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3c52644eff3f4a4fb8b6259daac7ff297a555421/e2e/25699400-7c3f-4264-a96c-bd7921a169d5.md", "", "", "836d1c83-52d1-4579-9d23-7f3bdff4659d.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f696022f6e3b5d387715a091f95e4dadf4e46478/e2e/da1b3d5f-85d4-4d96-bfee-43fbf7b22681.md", "", "", "ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md")
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3c52644eff3f4a4fb8b6259daac7ff297a555421/e2e/da1b3d5f-85d4-4d96-bfee-43fbf7b22681.md", "", "", "ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md")

# ===================== Sheet 3: de-de =====================
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.md"
$ws3.Range("G2").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.03b0573f4c3ad94c7c118a0f76853d51b4646da2.de-de.xlf"
$ws3.Range("H2").Value = "2016-09-01 09:20:15"
$ws3.Range("I2").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.md"
$ws3.Range("J2").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.03b0573f4c3ad94c7c118a0f76853d51b4646da2.de-de.xlf"
$ws3.Range("K2").Value = "2016-09-01 09:20:43"

$ws3.Range("A3").Value = "ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md"
$ws3.Range("G3").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.03b0573f4c3ad94c7c118a0f76853d51b4646da2.de-de.xlf"
$ws3.Range("H3").Value = "2016-09-01 09:20:15"
$ws3.Range("I3").Value = "ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md"
$ws3.Range("J3").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.03b0573f4c3ad94c7c118a0f76853d51b4646da2.de-de.xlf"
$ws3.Range("K3").Value = "2016-09-01 09:20:43"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f696022f6e3b5d387715a091f95e4dadf4e46478/e2e/25699400-7c3f-4264-a96c-bd7921a169d5.md", "", "", "836d1c83-52d1-4579-9d23-7f3bdff4659d.md")
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/af8d38cea75e92009d436a68111a43a0c5793587/e2e/25699400-7c3f-4264-a96c-bd7921a169d5.md", "", "", "836d1c83-52d1-4579-9d23-7f3bdff4659d.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f696022f6e3b5d387715a091f95e4dadf4e46478/e2e/da1b3d5f-85d4-4d96-bfee-43fbf7b22681.md", "", "", "ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md")
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/af8d38cea75e92009d436a68111a43a0c5793587/e2e/da1b3d5f-85d4-4d96-bfee-43fbf7b22681.md", "", "", "ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md")

Write-Host "Handback status report regenerated."
